# Actualización automática 2025-09-22 15:10:09
#
# A new salesperson, "PALMA PICO OSCAR FILIDEL", is inserted alphabetically
# (between MOROCHO PLAZA SHIRLEY AURELIA and PAREDES POVEDA TATIANA VERONICA)
# into both the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, as a new row
# 35. All of the following rows shift down by one, and the trailing totals
# row (which lived at row 44) moves to row 45. The totals text labels on the
# "VENTAS POR GRUPO" sheet ("X de 42") need to be updated to reflect the
# bigger population ("X de 43").

$wb = $excel.ActiveWorkbook

$xlShiftDown = -4121
$newName = "PALMA PICO OSCAR FILIDEL"
$office = "OFICINA-CATAECSA"

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" - columns A:R, data rows 2:44, totals row 44
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new blank row above the old row 35 (PAREDES POVEDA...), pushing
# every following row (including the totals row) down by one.
$ws1.Rows.Item(35).Insert($xlShiftDown)

# Populate the freshly inserted row 35 with the new salesperson; all of
# their sales figures are zero.
$ws1.Range("A35").Value = $office
$ws1.Range("B35").Value = $newName
$ws1.Range("C35:R35").Value = 0

# The summary row that used to be row 44 is now row 45; refresh its
# "X de 42" labels to "X de 43" now that there is one more salesperson.
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(45, $col)
    $cell.Value = $cell.Text.Replace("de 42", "de 43")
}

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL" - columns A:G, data rows 2:44, totals row 44
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(35).Insert($xlShiftDown)

$ws2.Range("A35").Value = $office
$ws2.Range("B35").Value = $newName
$ws2.Range("C35:G35").Value = 0

$wb.Save()
